$d = $word.ActiveDocument

# Fill in the contact-details table row for Timothy (name, student number,
# e-mail hyperlink, phone number).
$tbl = $d.Tables(1)
$row = $tbl.Rows(4)

$cellName = $row.Cells(1).Range
$cellName.Collapse(0)  # wdCollapseEnd
$cellName.MoveEnd(1, -1)
$cellName.InsertAfter("Benschop")

$cellNum = $row.Cells(2).Range
$cellNum.MoveEnd(1, -1)
$cellNum.Text = "18067301"

$cellEmail = $row.Cells(3).Range
$cellEmail.MoveEnd(1, -1)
$d.Hyperlinks.Add($cellEmail, "mailto:18067301@student.hhs.nl", $null, $null, "18067301@student.hhs.nl") | Out-Null

$cellPhone = $row.Cells(4).Range
$cellPhone.MoveEnd(1, -1)
$cellPhone.Text = "0631521230"
